$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.769.43'
$ws.Range('E2').Value = '  -2.75%  '
$ws.Range('D3').Value = '1.747.59'
$ws.Range('E3').Value = '  -5.23%  '
$ws.Range('D4').Value = '''1.004'
$ws.Range('E4').Value = '  +0.24%  '
$ws.Range('D5').Value = '''238.45'
$ws.Range('E5').Value = '  -8.87%  '
$ws.Range('D6').Value = '''1.005'
$ws.Range('E6').Value = '  +0.34%  '
$ws.Range('D7').Value = '''0.5013'
$ws.Range('E7').Value = '  -6.04%  '
$ws.Range('D8').Value = '''41.73'
$ws.Range('E8').Value = '  -6.92%  '
$ws.Range('D9').Value = '''0.2632'
$ws.Range('E9').Value = '  -13.39%  '
$ws.Range('D10').Value = '''0.06129'
$ws.Range('E10').Value = '  -11.10%  '
$ws.Range('D11').Value = '1.750.41'
$ws.Range('E11').Value = '  -6.40%  '
$ws.Range('D12').Value = '''0.06954'
$ws.Range('E12').Value = '  -9.43%  '
$ws.Range('D13').Value = '''15.14'
$ws.Range('E13').Value = '  -16.47%  '
$ws.Range('D14').Value = '''4.489'
$ws.Range('E14').Value = '  -10.16%  '
$ws.Range('B15').Value = 'Polygon'
$ws.Range('C15').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D15').Value = '''0.5872'
$ws.Range('E15').Value = '  -20.64%  '
$ws.Range('B16').Value = 'Litecoin'
$ws.Range('C16').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D16').Value = '''76.57'
$ws.Range('E16').Value = '  -14.72%  '
$ws.Range('D17').Value = '''1.004'
$ws.Range('E17').Value = '  +0.24%  '
$ws.Range('D18').Value = '''1.002'
$ws.Range('E18').Value = '  +0.03%  '
$ws.Range('D19').Value = '25.791.94'
$ws.Range('E19').Value = '  -2.75%  '
$ws.Range('D20').Value = '''11.64'
$ws.Range('E20').Value = '  -16.73%  '
$ws.Range('D21').Value = '''0.000006740'
$ws.Range('E21').Value = '  -15.22%  '
$ws.Range('D22').Value = '1.971.90'
$ws.Range('E22').Value = '  -5.62%  '
$ws.Range('D23').Value = '''4.063'
$ws.Range('E23').Value = '  -12.01%  '
$ws.Range('D24').Value = '''8.089'
$ws.Range('E24').Value = '  -13.10%  '
$ws.Range('D25').Value = '''5.099'
$ws.Range('E25').Value = '  -14.92%  '
$ws.Range('D26').Value = '''137.99'
$ws.Range('E26').Value = '  -3.71%  '
$ws.Range('D27').Value = '''1.536'
$ws.Range('E27').Value = '  -8.94%  '
$ws.Range('D28').Value = '''1.826'
$ws.Range('E28').Value = '  -17.07%  '
$ws.Range('D29').Value = '''14.86'
$ws.Range('E29').Value = '  -12.54%  '
$ws.Range('D30').Value = '''102.84'
$ws.Range('E30').Value = '  -7.19%  '
$ws.Range('D31').Value = '''0.08112'
$ws.Range('E31').Value = '  -7.73%  '
$ws.Range('D32').Value = '''3.741'
$ws.Range('E32').Value = '  -12.23%  '
$ws.Range('D33').Value = '''3.458'
$ws.Range('E33').Value = '  -14.91%  '
$ws.Range('D34').Value = '''0.04485'
$ws.Range('E34').Value = '  -6.63%  '
$ws.Range('D35').Value = '''1.004'
$ws.Range('E35').Value = '  +0.31%  '
$ws.Range('D36').Value = '''2.653'
$ws.Range('E36').Value = '  -9.50%  '
$ws.Range('D37').Value = '''0.9861'
$ws.Range('E37').Value = '  -13.11%  '
$ws.Range('D38').Value = '''0.6042'
$ws.Range('E38').Value = '  -16.94%  '
$ws.Range('D39').Value = '''2.646'
$ws.Range('E39').Value = '  -14.92%  '
$ws.Range('D40').Value = '''0.01546'
$ws.Range('E40').Value = '  -10.09%  '
$ws.Range('B41').Value = 'PaxDollar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D41').Value = '''1.006'
$ws.Range('E41').Value = '  +0.47%  '
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D42').Value = '''1.914'
$ws.Range('E42').Value = '  -16.99%  '
$ws.Range('D43').Value = '''104.17'
$ws.Range('E43').Value = '  -3.58%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').Value = '''5.160'
$ws.Range('E44').Value = '  -12.29%  '
$ws.Range('B45').Value = 'TheSandbox'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D45').Value = '''0.3791'
$ws.Range('E45').Value = '  -20.60%  '
$ws.Range('D46').Value = '''0.7239'
$ws.Range('E46').Value = '  -21.04%  '
$ws.Range('D47').Value = '''0.05299'
$ws.Range('E47').Value = '  -8.61%  '
$ws.Range('D48').Value = '''0.1105'
$ws.Range('E48').Value = '  -11.05%  '
$ws.Range('B49').Value = 'Elrond'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D49').Value = '''30.11'
$ws.Range('E49').Value = '  -13.58%  '
$ws.Range('B50').Value = 'Aptos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D50').Value = '''5.909'
$ws.Range('E50').Value = '  -21.04%  '
$ws.Range('D51').Value = '''52.31'
$ws.Range('E51').Value = '  -13.48%  '
